$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A first (comments), in row order
$ws.Range("A1").Value = "i hate this product"
$ws.Range("A2").Value = " this product is so fantastic its great"
$ws.Range("A3").Value = " i will never use this product again"

# Column B next (sentiment labels), in row order
$ws.Range("B1").Value = "Very Negative"
$ws.Range("B2").Value = "Very Positive"
$ws.Range("B3").Value = "Very Negative"

# Column C last (scores)
$ws.Range("C1").Value = 1
$ws.Range("C2").Value = 0.98
$ws.Range("C3").Value = 0.5
